$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 4-6: Privacy(4) -> row6, AI literacy(2) -> row4, Transparency(2) -> row5
$ws.Range("A4").Value = "AI literacy"
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = "Transparency"
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = "Privacy"
$ws.Range("C6").Value = 4
